# Fill in the "Name", "Roll Number" and "Branch worked On" cells for the
# Team Lead row (row 4) of the Group Member Details table with the
# author's own details, matching the formatting of the surrounding cells.

function Set-CellRun {
    param($row, $col, $paraId, $textId, $rsidR, $rsidRPr, $rsidRDefault, $rsidP, $runRPr, $text)

    $cellRange = $word.ActiveDocument.Tables(1).Cell($row, $col).Range

    $ooxml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' +
        '<w:p w14:paraId="' + $paraId + '" w14:textId="' + $textId + '" w:rsidR="' + $rsidR + '" w:rsidRPr="' + $rsidRPr + '" w:rsidRDefault="' + $rsidRDefault + '" w:rsidP="' + $rsidP + '">' +
        '<w:pPr><w:rPr>' + $runRPr + '</w:rPr></w:pPr>' +
        '<w:r><w:rPr>' + $runRPr + '</w:rPr><w:t>' + $text + '</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    [void]$cellRange.InsertXML($ooxml)
}

$boldCalibriRPr = '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/>'
$calibriRPr = '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>'

# Name cell
Set-CellRun 4 2 "4B0C437D" "395E9462" "00250253" "00B84C04" "00250253" "004F6518" $boldCalibriRPr "Muhammad Adil Azan"

# Roll Number cell
Set-CellRun 4 3 "1BC8304E" "77777777" "00250253" "00B84C04" "00250253" "004F6518" $boldCalibriRPr "L1F20BSSE0592"

# Branch worked On cell
Set-CellRun 4 5 "7809FB32" "77777777" "00250253" "002D05DF" "00250253" "004F6518" $calibriRPr "develop"
